$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Development Log entry: copy the formatting from the row above (row 11)
# so the date/time/number formats & borders match the rest of the table,
# then fill in this session's values.
$ws.Range("C11:G11").Copy()
$ws.Range("C15:G15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C15").Value = 44637                    # 17/03/2022
$ws.Range("D15").Value = 0.52083333333333337      # 12:30
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "Implementation"
$ws.Range("G15").Value = "Implementation"

# Update the saved active-cell selection to H15
$ws.Range("H15").Select()
